# Applies revisions to the UR Cover Sheets workbook:
#  - renames the five unit-number tabs
#  - updates inshopped/dead dates, unit numbers, due flags and defect notes
#    on each of those tabs

$wb = $excel.ActiveWorkbook

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting numeric- or date-looking strings into numbers/dates.
# A leading apostrophe forces "stored as text" semantics, matching the
# inlineStr cells expected by the target workbook.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

$sheets = @(
    @{ Index = 2; OldName = "1050"; NewName = "9905" },
    @{ Index = 3; OldName = "3651"; NewName = "5673" },
    @{ Index = 4; OldName = "3657"; NewName = "7030" },
    @{ Index = 5; OldName = "5178"; NewName = "3337" },
    @{ Index = 6; OldName = "7534"; NewName = "4023" }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Index)
    $ws.Name = $s.NewName
}

# ---- Sheet "9905" (was "1050") ----
$ws = $wb.Worksheets.Item(2)
Set-TextValue $ws.Range("F1") "01-16-2020"
Set-TextValue $ws.Range("A2") "9905"
Set-TextValue $ws.Range("C5") "06-05-2020"
$ws.Range("B23").Value = "FUEL LEAK"
$ws.Range("B24").Value = "Fuel leaking right side"

# ---- Sheet "5673" (was "3651") ----
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws.Range("F1") "01-16-2020"
Set-TextValue $ws.Range("A2") "5673"
$ws.Range("F3").Value = "Y"
Set-TextValue $ws.Range("C5") "03-19-2020"
$ws.Range("B23").Value = "AIR BRAKE"
$ws.Range("B24").Value = "Automatic leaks"
$ws.Range("B26").Value = ""
$ws.Range("B27").Value = ""

# ---- Sheet "7030" (was "3657") ----
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("F1") "01-16-2020"
Set-TextValue $ws.Range("A2") "7030"
$ws.Range("F3").Value = "Y"
Set-TextValue $ws.Range("C5") "05-18-2020"
$ws.Range("B23").Value = "Check for tags"
$ws.Range("B24").Value = ""

# ---- Sheet "3337" (was "5178") ----
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("F1") "01-16-2020"
Set-TextValue $ws.Range("A2") "3337"
Set-TextValue $ws.Range("C5") "02-16-2020"
$ws.Range("F6").Value = "Y"
$ws.Range("B23").Value = "TURBO"
$ws.Range("B24").Value = "Grinding when turbo breaker is closed"
$ws.Range("B26").Value = "Carbody"
$ws.Range("B27").Value = "Rear truck dynamic cover open"

# ---- Sheet "4023" (was "7534") ----
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("F1") "01-16-2020"
Set-TextValue $ws.Range("A2") "4023"
Set-TextValue $ws.Range("C5") "06-02-2020"
$ws.Range("B23").Value = "STRUCK TREE"
$ws.Range("B24").Value = ""
